$wb = $excel.ActiveWorkbook

# Rename the "Include from ..." sheets to "Include #N"
$wb.Worksheets.Item("Include from UsageContextType").Name = "Include #0"
$wb.Worksheets.Item("Include from TRE-R67-TypeStru").Name = "Include #1"
$wb.Worksheets.Item("Include from TRE-R04-TypeSavo").Name = "Include #2"
$wb.Worksheets.Item("Include from TRE-R288-TypePro").Name = "Include #3"

# Update the Date value on the Metadata sheet
$ws = $wb.Worksheets.Item("Metadata")
$ws.Range("B8").Value = "2024-09-13T14:28:16+00:00"
